# "integrate database all in one"
#
# The balance-sheet sheet keeps a rolling window of 5 fiscal periods in
# columns D:H (column C is a blank spacer next to the row labels in B).
# This update rolls the window forward one period:
#   - the oldest period (12 ماهه منتهی به 1397/06, column D) is dropped
#   - the former column E period (1398/05) slides left into column D
#   - a brand new period (1399/05) is inserted as the new column E,
#     together with its own freshly reported figures
#   - columns F, G, H (1399/12, 1400/08, 1401/08) are untouched
#
# Net effect on every data row: D <- old E, E <- newly reported value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: period headers ---------------------------------------------
$ws.Range("D8").Value = "12 ماهه منتهی به 1398/05"
$ws.Range("E8").Value = "12 ماهه منتهی به 1399/05"

# --- Row 9: publish dates ------------------------------------------------
$ws.Range("D9").Value = "1399-07-24 (8)"
$ws.Range("E9").Value = "1400-01-30 (5)"
$ws.Range("F9").Value = "1400-10-08 (6)"
$ws.Range("H9").Value = "1402-03-30 (5)"

# --- دارایی (Assets) ------------------------------------------------------
# موجودی نقد
$ws.Range("D12").Value = 238048
$ws.Range("E12").Value = 527797
# سرمایه گذاری کوتاه مدت
$ws.Range("D13").Value = 102212
$ws.Range("E13").Value = 151529
# دریافتنی‌های تجاری و سایر دریافتنی‌ها
$ws.Range("D14").Value = 302775
$ws.Range("E14").Value = 253631
# موجودی مواد و کالا
$ws.Range("D15").Value = 237811
$ws.Range("E15").Value = 480708
# پیش پرداخت ها
$ws.Range("D16").Value = 94713
$ws.Range("E16").Value = 36614
# دارایی های نگهداری شده برای فروش
$ws.Range("D17").Value = 5113
$ws.Range("E17").Value = 0
# جمع داراییهای جاری
$ws.Range("D18").Value = 980672
$ws.Range("E18").Value = 1450279
# حسابها و اسناد دریافتنی تجاری بلند مدت (E already 0, unchanged)
$ws.Range("D19").Value = 0
# سرمایه گذاریهای بلند مدت
$ws.Range("D20").Value = 135400
$ws.Range("E20").Value = 170220
# داراییهای ثابت مشهود
$ws.Range("D22").Value = 837329
$ws.Range("E22").Value = 792828
# داراییهای نامشهود
$ws.Range("D23").Value = 11661
$ws.Range("E23").Value = 42668
# جمع داراییهای غیرجاری
$ws.Range("D26").Value = 984390
$ws.Range("E26").Value = 1005716
# جمع داراییها
$ws.Range("D27").Value = 1965062
$ws.Range("E27").Value = 2455995

# --- بدهی (Liabilities) ---------------------------------------------------
# پرداختنی‌های تجاری و سایر پرداختنی‌ها
$ws.Range("D29").Value = 186237
$ws.Range("E29").Value = 264137
# پیش دریافتها
$ws.Range("D31").Value = 14528
$ws.Range("E31").Value = 32487
# ذخیره مالیات بر درامد
$ws.Range("D32").Value = 71156
$ws.Range("E32").Value = 69319
# سود سهام پیشنهادی و پرداختنی
$ws.Range("D33").Value = 11041
$ws.Range("E33").Value = 25533
# حصه جاری تسهیلات مالی دریافتی
$ws.Range("D34").Value = 327890
$ws.Range("E34").Value = 321172
# جمع بدهیهای جاری
$ws.Range("D37").Value = 610852
$ws.Range("E37").Value = 712648
# تسهیلات مالی دریافتی بلند مدت
$ws.Range("D40").Value = 135025
$ws.Range("E40").Value = 96303
# ذخیره مزایای پایان خدمت
$ws.Range("D41").Value = 17927
$ws.Range("E41").Value = 45388
# جمع بدهیهای غیر جاری
$ws.Range("D42").Value = 152952
$ws.Range("E42").Value = 141691
# جمع بدهیهای جاری و غیر جاری
$ws.Range("D43").Value = 763804
$ws.Range("E43").Value = 854339

# --- حقوق صاحبان سهام (Equity) -------------------------------------------
# اندوخته قانونی (E already 50000, unchanged)
$ws.Range("D50").Value = 50000
# سود (زیان) انباشته
$ws.Range("D56").Value = 651258
$ws.Range("E56").Value = 1051656
# جمع حقوق صاحبان سهام
$ws.Range("D57").Value = 1201258
$ws.Range("E57").Value = 1601656
# جمع بدهیها و حقوق صاحبان سهام
$ws.Range("D58").Value = 1965062
$ws.Range("E58").Value = 2455995
